$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row 4 with the new mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(4, 1).Value = "Demo inplannen"
$logs.Cells.Item(4, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item(4, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item(4, 6).Value = "2025-08-19 19:29:14"
$logs.Cells.Item(4, 7).Value = "Nee"
$logs.Cells.Item(4, 8).Value = "Ja"
$logs.Cells.Item(4, 9).Value = "Nee"
$logs.Cells.Item(4, 10).Value = "Nee"

# Extend the conditional formatting ranges (D/G/H/I/J) down to the new row 4
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D4"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G4"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H4"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I4"))
$logs.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J4"))

# --- Sheet "Dashboard": bump the count for this category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 3
